$wb = $excel.ActiveWorkbook

# --- Sheet 1: GNG_TO ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "GNG_TO-16509961675314205"
$ws1.Range("B2").Value = "go_stims-16509961674993823.csv"
$ws1.Range("B3").Value = "GNG_stims-16509961675154276.csv"
$ws1.Range("B4").Value = "go_stims-16509961675154276.csv"
$ws1.Range("B5").Value = "GNG_stims-16509961675314205.csv"

# --- Sheet 2: NB_TO ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "NB_TO-16509961709153838"
$ws2.Range("B2").Value = "OB-16509961694833794.csv"
$ws2.Range("B3").Value = "TB-16509961704034219.csv"
$ws2.Range("B4").Value = "ZB-match_4-16509961677394195.csv"
$ws2.Range("B5").Value = "TB-16509961708913774.csv"
$ws2.Range("B6").Value = "TB-16509961699074192.csv"
$ws2.Range("B7").Value = "ZB-match_0-16509961675633843.csv"
$ws2.Range("B8").Value = "ZB-match_7-16509961677633908.csv"
$ws2.Range("B9").Value = "OB-16509961685233817.csv"
$ws2.Range("B10").Value = "OB-16509961680353801.csv"

# --- Sheet 3: RS_TO ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "RS_TO-16509961709153838"
$ws3.Range("B2").Value = "eyes closed"
$ws3.Range("B3").Value = "eyes open"

# --- Sheet 4: TOL_TO ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Name = "TOL_TO-16509961709874165"
$ws4.Range("B2").Value = "MM_stims-16509961709474149.csv"
$ws4.Range("B3").Value = "ZM_stims-16509961709233828.csv"
$ws4.Range("B4").Value = "MM_stims-16509961709713814.csv"
$ws4.Range("B5").Value = "ZM_stims-16509961709474149.csv"
$ws4.Range("B6").Value = "MM_stims-16509961709874165.csv"
$ws4.Range("B7").Value = "ZM_stims-16509961709713814.csv"

# --- Sheet 5: vSAT_TO ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Name = "vSAT_TO-16509961710514162"
$ws5.Range("B2").Value = "vSAT_stims-16509961710354185.csv"
$ws5.Range("B3").Value = "SAT_stims-1650996171003415.csv"
$ws5.Range("B4").Value = "SAT_stims-16509961709874165.csv"
$ws5.Range("B5").Value = "vSAT_stims-16509961710194142.csv"
